# KIBON-2763: add column "Status" to report "Gemeinden" (sheet "Angaben pro Periode")

$wb = $excel.ActiveWorkbook
$wsPeriode = $wb.Worksheets.Item("Angaben pro Periode")

# Insert a new column before the current column E ("limitierungKita"), shifting
# the existing E..K columns one to the right (F..L). Insert() copies the
# neighbouring column's formatting onto the new column/cells, matching what
# Excel does.
$wsPeriode.Columns.Item(5).Insert()

# New header cell (row 5) and example-value cell (row 6) for the new "Status" column.
$wsPeriode.Cells.Item(5, 5).Value = "{gemeindeKennzahlenStatusTitle}"
$wsPeriode.Cells.Item(6, 5).Value = "{gemeindeKennzahlenStatus}"

# The new column is not a "best fit" column like its neighbours - it gets an
# explicit custom width instead.
$wsPeriode.Columns.Item(5).ColumnWidth = 19

# The last column (L6, previously K6) loses its (unused/ghost) formatting.
$wsPeriode.Cells.Item(6, 12).ClearFormats()

# Update the remembered selection on this sheet.
$wsPeriode.Range("E17").Select()

# "Angaben pro Periode" becomes the active tab/sheet of the workbook.
$wsPeriode.Activate()
